$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - keep name, update values
$ws.Range("B3").Value = 0.9930717586779605
$ws.Range("C3").Value = 0.9934125538928393
$ws.Range("D3").Value = 0.9851299774447799

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9914104756139079
$ws.Range("C4").Value = 0.9917526044288346
$ws.Range("D4").Value = 0.9882868978950633

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9962295202574335
$ws.Range("C5").Value = 0.9956081812058916
$ws.Range("D5").Value = 0.9953257782265746
